# Apply cryptocurrency price/volume updates to Sheet1 (cryptos.xlsx)
# Source: scheduled GitHub Actions data refresh commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.713.16"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.878.55"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3962"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.02"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08031"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.024"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "1.905.03"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.967"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.174"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06615"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "27.721.87"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.506"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "2.124.67"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.66"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9717"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09553"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.453"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.628"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.305"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06118"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02264"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.234"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.170"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6004"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1900"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.250"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5694"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.20"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.398"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.934"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06820"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000314"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.70%  "

Write-Output "Applied 95 cell updates to Sheet1"
